$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 currently holds the "Knärot" (Goodyera repens) record and row 16
# holds the "Bronshjon" (Callidium coriaceum) record. This edit swaps the
# two records' data (id numbers, names, coordinates, ...) while keeping
# them on the same row numbers 15/16 within the sheet.

# --- capture the simple scalar values that exist on both rows (these get
#     swapped in place) ---
$A15 = $ws.Range("A15").Value()
$B15 = $ws.Range("B15").Value()
$D15 = $ws.Range("D15").Value()
$E15 = $ws.Range("E15").Value()
$F15 = $ws.Range("F15").Value()
$G15 = $ws.Range("G15").Value()
$H15 = $ws.Range("H15").Value()
$Q15 = $ws.Range("Q15").Value()
$R15 = $ws.Range("R15").Value()

$A16 = $ws.Range("A16").Value()
$B16 = $ws.Range("B16").Value()
$D16 = $ws.Range("D16").Value()
$E16 = $ws.Range("E16").Value()
$F16 = $ws.Range("F16").Value()
$G16 = $ws.Range("G16").Value()
$H16 = $ws.Range("H16").Value()
$Q16 = $ws.Range("Q16").Value()
$R16 = $ws.Range("R16").Value()

$ws.Range("A15").Value = $A16
$ws.Range("B15").Value = $B16
$ws.Range("D15").Value = $D16
$ws.Range("E15").Value = $E16
$ws.Range("F15").Value = $F16
$ws.Range("G15").Value = $G16
$ws.Range("H15").Value = $H16
$ws.Range("Q15").Value = $Q16
$ws.Range("R15").Value = $R16

$ws.Range("A16").Value = $A15
$ws.Range("B16").Value = $B15
$ws.Range("D16").Value = $D15
$ws.Range("E16").Value = $E15
$ws.Range("F16").Value = $F15
$ws.Range("G16").Value = $G15
$ws.Range("H16").Value = $H15
$ws.Range("Q16").Value = $Q15
$ws.Range("R16").Value = $R15

# --- the "Bronshjon" record additionally carries values in J/K/L/M/N/AF
#     (Enhet/Ålder-Stadium/Kön/Aktivitet/Metod/Bestämningsmetod) that the
#     "Knärot" record does not have at all. Move those cells (currently on
#     row 16) over to row 15 with Cut, which relocates them (and leaves the
#     row 16 side blank / without those cells), matching the swap of the
#     two records.
$ws.Range("J16").Cut($ws.Range("J15"))
$ws.Range("K16").Cut($ws.Range("K15"))
$ws.Range("L16").Cut($ws.Range("L15"))
$ws.Range("M16").Cut($ws.Range("M15"))
$ws.Range("N16").Cut($ws.Range("N15"))
$ws.Range("AF16").Cut($ws.Range("AF15"))
